$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared/rich-text string edits ---
# "Volume 33   Number  5" -> "...6" (A8)
$ws.Range("A8").Characters(21, 1).Text = "6"

# "Report Covering the Week  1/26/2026  Through  2/1/2026" -> updated dates (C9)
$ws.Range("C9").Characters(27, 9).Text = "2/2/2026"
$ws.Range("C9").Characters(46, 8).Text = "2/8/2026"

# --- Crime-statistics table value/type updates (rows 15-28) ---
$ws.Range("F15").NumberFormat = '@'
$ws.Range("F15").Value = '0'
$ws.Range("L15").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("L15").Value = 0
$ws.Range("C16").NumberFormat = '#,##0'
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = -71.428571428571
$ws.Range("I16").Value = 6
$ws.Range("J16").Value = 11
$ws.Range("K16").Value = -45.454545454545
$ws.Range("L16").Value = -53.846153846153
$ws.Range("M16").Value = -71.428571428571
$ws.Range("N16").Value = -90.625
$ws.Range("C17").Value = 2
$ws.Range("D17").NumberFormat = '#,##0'
$ws.Range("D17").Value = 5
$ws.Range("E17").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = -33.333333333333
$ws.Range("I17").Value = 11
$ws.Range("J17").Value = 11
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = -8.333333333333
$ws.Range("M17").Value = 83.333333333333
$ws.Range("N17").Value = -52.173913043478
$ws.Range("C18").NumberFormat = '#,##0'
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -78.571428571428
$ws.Range("I18").Value = 5
$ws.Range("J18").Value = 18
$ws.Range("K18").Value = -72.222222222222
$ws.Range("M18").Value = -77.272727272727
$ws.Range("N18").Value = -96.478873239436
$ws.Range("C19").Value = 7
$ws.Range("E19").Value = -22.222222222222
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 37
$ws.Range("H19").Value = -16.216216216216
$ws.Range("I19").Value = 42
$ws.Range("J19").Value = 57
$ws.Range("K19").Value = -26.315789473684
$ws.Range("L19").Value = -37.31343283582
$ws.Range("M19").Value = 75
$ws.Range("N19").Value = 0
$ws.Range("C20").Value = 1
$ws.Range("D20").NumberFormat = '#,##0'
$ws.Range("D20").Value = 2
$ws.Range("E20").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E20").Value = -50
$ws.Range("I20").Value = 10
$ws.Range("J20").Value = 5
$ws.Range("K20").Value = 100
$ws.Range("L20").Value = 100
$ws.Range("M20").Value = 25
$ws.Range("N20").Value = -88.636363636363
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = -36.842105263157
$ws.Range("F21").Value = 49
$ws.Range("G21").Value = 69
$ws.Range("H21").Value = -28.985507246376
$ws.Range("I21").Value = 75
$ws.Range("J21").Value = 102
$ws.Range("K21").Value = -26.470588235294
$ws.Range("L21").Value = -36.440677966101
$ws.Range("M21").Value = -7.407407407407
$ws.Range("N21").Value = -79.224376731301
$ws.Range("C22").NumberFormat = '@'
$ws.Range("C22").Value = '0'
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 6
$ws.Range("K22").Value = -66.666666666666
$ws.Range("F23").NumberFormat = '@'
$ws.Range("F23").Value = '0'
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = -100
$ws.Range("J23").Value = 4
$ws.Range("K23").Value = -75
$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = -30
$ws.Range("F24").Value = 39
$ws.Range("G24").Value = 63
$ws.Range("H24").Value = -38.095238095238
$ws.Range("I24").Value = 62
$ws.Range("J24").Value = 84
$ws.Range("K24").Value = -26.190476190476
$ws.Range("L24").Value = -37.373737373737
$ws.Range("M24").Value = 31.914893617021
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -80
$ws.Range("F25").Value = 12
$ws.Range("G25").Value = 41
$ws.Range("H25").Value = -70.731707317073
$ws.Range("I25").Value = 23
$ws.Range("J25").Value = 60
$ws.Range("K25").Value = -61.666666666666
$ws.Range("L25").Value = -51.063829787234
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = 33.333333333333
$ws.Range("F26").Value = 21
$ws.Range("G26").Value = 16
$ws.Range("H26").Value = 31.25
$ws.Range("I26").Value = 32
$ws.Range("J26").Value = 22
$ws.Range("K26").Value = 45.454545454545
$ws.Range("L26").Value = -3.030303030303
$ws.Range("M26").Value = 77.777777777777
$ws.Range("F27").NumberFormat = '@'
$ws.Range("F27").Value = '0'
$ws.Range("L27").Value = -50
$ws.Range("D28").NumberFormat = '@'
$ws.Range("D28").Value = '0'
$ws.Range("E28").NumberFormat = '@'
$ws.Range("E28").Value = '***.*'
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -75
